$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Remove Sheet2 and Sheet3, leaving only Sheet1 ---
$null = $wb.Worksheets.Item("Sheet2").Delete()
$null = $wb.Worksheets.Item("Sheet3").Delete()

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()

# --- Fix up the query string referenced by column A (shared across A2:A11) ---
# This also ripples into the cached values of the G2:G11 "INSERT INTO ..." formulas.
for ($r = 2; $r -le 11; $r++) {
    $ws1.Cells.Item($r, 1).Value = "[VolTeer].[Vol].[tblVolunteer]"
}

# --- New trailing blank row (row 15), matching the existing blank row 12's height ---
$ws1.Rows.Item(15).RowHeight = 13.8

# --- Update the workbook-level tab ratio (bookViews/workbookView@tabRatio) ---
$excel.Windows.Item(1).TabRatio = 0.917

# --- Update sheet view: scroll so column D is left-most, and select E3 ---
$null = $excel.Goto($ws1.Range("D1"), $true)
$null = $ws1.Range("E3").Select()
